$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: Real_Median -> -1 for rows 2-9
$ws.Range("B2:B9").Value = -1

# Column C: Real_IQR -> 0 for rows 2-9
$ws.Range("C2:C9").Value = 0

# Column D: Pred_Median -> new values
$ws.Range("D2").Value = 0.00005927681922912598
$ws.Range("D3").Value = 0.00005584955215454102
$ws.Range("D4").Value = 0.00005315244197845459
$ws.Range("D5").Value = 0.00006470084190368652
$ws.Range("D6").Value = 0.00005459785461425781
$ws.Range("D7").Value = 0.00005787611007690430
$ws.Range("D8").Value = 0.00005516409873962402
$ws.Range("D9").Value = 0.00005093216896057129

# Column E: Pred_IQR -> new values
$ws.Range("E2").Value = 0.00009092688560485840
$ws.Range("E3").Value = 0.00010088086128234861
$ws.Range("E4").Value = 0.00010339170694351200
$ws.Range("E5").Value = 0.00009348988533020020
$ws.Range("E6").Value = 0.00008073449134826660
$ws.Range("E7").Value = 0.00010156631469726560
$ws.Range("E8").Value = 0.00008788704872131348
$ws.Range("E9").Value = 0.00008669495582580566
